$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Peer  and self assessment")

# Fill in Daniel's assessment for Criterion 1 Online collaboration (row 7)
$ws.Range("B7").Value = "Sufficient"
$ws.Range("C7").Value = "Not very active on discord (currently 7 posts). Missed a meeting."

# Fill in Daniel's assessment for Criterion 1 International Collaboration (row 20)
$ws.Range("B20").Value = "Excellent"
$ws.Range("C20").Value = "Active collaborator, motivated"

# Update the active selection to reflect the last edited cell
$ws.Activate()
$ws.Range("C7").Select()
